# Updates the cryptos list (Price / Volume(1h) columns) as produced by the
# scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> @{ D = new price text (or $null if unchanged); E = new volume text }
# "@"  -> the Price value parses as a plain number in Excel, so it must be
#         forced to Text format first or the trailing/insignificant zeros
#         (e.g. "533.20", "1.00") would be silently dropped.
$updates = @(
    @{ Row = 2;  D = "71.691.45";  DFormat = $null; E = "  +4.79%  " }
    @{ Row = 3;  D = "4.050.18";   DFormat = $null; E = "  +5.10%  " }
    @{ Row = 4;  D = $null;        DFormat = $null; E = "  +0.04%  " }
    @{ Row = 5;  D = "533.20";     DFormat = "@";   E = "  +2.23%  " }
    @{ Row = 6;  D = "153.51";     DFormat = "@";   E = "  +9.42%  " }
    @{ Row = 7;  D = "0.694";      DFormat = "@";   E = "  +14.28%  " }
    @{ Row = 8;  D = $null;        DFormat = $null; E = "  +0.04%  " }
    @{ Row = 9;  D = $null;        DFormat = $null; E = "  +7.22%  " }
    @{ Row = 10; D = "0.176";      DFormat = "@";   E = "  +6.04%  " }
    @{ Row = 11; D = $null;        DFormat = $null; E = "  +5.07%  " }
    @{ Row = 12; D = "49.06";      DFormat = "@";   E = "  +18.71%  " }
    @{ Row = 13; D = $null;        DFormat = $null; E = "  +6.82%  " }
    @{ Row = 14; D = "4.691.52";   DFormat = $null; E = "  +4.94%  " }
    @{ Row = 15; D = "4.042.72";   DFormat = $null; E = "  +4.59%  " }
    @{ Row = 16; D = "14.51";      DFormat = "@";   E = "  +3.16%  " }
    @{ Row = 17; D = "21.08";      DFormat = "@";   E = "  -1.27%  " }
    @{ Row = 18; D = $null;        DFormat = $null; E = "  +2.49%  " }
    @{ Row = 19; D = $null;        DFormat = $null; E = "  +0.10%  " }
    @{ Row = 20; D = "71.813.81";  DFormat = $null; E = "  +4.84%  " }
    @{ Row = 21; D = "437.02";     DFormat = "@";   E = "  +5.38%  " }
    @{ Row = 22; D = "3.74";       DFormat = "@";   E = "  +7.79%  " }
    @{ Row = 23; D = "99.47";      DFormat = "@";   E = "  +14.90%  " }
    @{ Row = 24; D = "14.84";      DFormat = "@";   E = "  +6.29%  " }
    @{ Row = 25; D = "4.23";       DFormat = "@";   E = "  +6.55%  " }
    @{ Row = 26; D = "11.44";      DFormat = "@";   E = "  -1.64%  " }
    @{ Row = 27; D = "10.92";      DFormat = "@";   E = "  +4.07%  " }
    @{ Row = 28; D = "37.34";      DFormat = "@";   E = "  +5.82%  " }
    @{ Row = 29; D = "5.83";       DFormat = "@";   E = "  +3.00%  " }
    @{ Row = 30; D = "3.55";       DFormat = "@";   E = "  +28.06%  " }
    @{ Row = 31; D = $null;        DFormat = $null; E = "  +3.93%  " }
    @{ Row = 32; D = $null;        DFormat = $null; E = "  +6.56%  " }
    @{ Row = 33; D = "679.59";     DFormat = "@";   E = "  +0.38%  " }
    @{ Row = 34; D = "6.77";       DFormat = "@";   E = "  +1.19%  " }
    @{ Row = 35; D = "66.44";      DFormat = "@";   E = "  +1.73%  " }
    @{ Row = 36; D = "43.11";      DFormat = "@";   E = "  +9.07%  " }
    @{ Row = 37; D = $null;        DFormat = $null; E = "  -3.66%  " }
    @{ Row = 38; D = "0.159";      DFormat = "@";   E = "  +6.85%  " }
    @{ Row = 39; D = "0.0₃0863";   DFormat = $null; E = "  +4.57%  " }
    @{ Row = 40; D = "3.44";       DFormat = "@";   E = "  -2.39%  " }
    @{ Row = 41; D = "1.00";       DFormat = "@";   E = "  +0.12%  " }
    @{ Row = 42; D = $null;        DFormat = $null; E = "  +5.65%  " }
    @{ Row = 43; D = $null;        DFormat = $null; E = "  -0.15%  " }
    @{ Row = 44; D = "3.17";       DFormat = "@";   E = "  +2.42%  " }
    @{ Row = 45; D = $null;        DFormat = $null; E = "  +9.15%  " }
    @{ Row = 46; D = $null;        DFormat = $null; E = "  -1.93%  " }
    @{ Row = 47; D = "3.42";       DFormat = "@";   E = "  +0.72%  " }
    @{ Row = 48; D = "9.58";       DFormat = "@";   E = "  +12.44%  " }
    @{ Row = 49; D = "3.06";       DFormat = "@";   E = "  +2.62%  " }
    @{ Row = 50; D = "3.38";       DFormat = "@";   E = "  +3.95%  " }
    @{ Row = 51; D = "0.000275";   DFormat = "@";   E = "  +3.05%  " }
)

foreach ($entry in $updates) {
    $row = $entry.Row
    if ($null -ne $entry.D) {
        $cell = $ws.Cells.Item($row, 4)
        if ($null -ne $entry.DFormat) {
            $cell.NumberFormat = $entry.DFormat
        }
        $cell.Value = $entry.D
    }
    $ws.Cells.Item($row, 5).Value = $entry.E
}
